# Natmi following Dr Hou advice
# Rebuild the Sema3a-Nrp2 LR-pairs sheet to include the "ECs" sending/target
# cluster alongside the existing "FAPs" and "sCs" clusters, expanding the
# table from 2x3 (6 data rows) to 3x3 (9 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Full data for the rebuilt table (rows 2..10) ----
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster
# E..T numeric metrics
$rows = @(
    @("ECs","Sema3a","Nrp2","ECs", 2,0.6666666666666666,1.397441,4.192323,0.6676161521996591,0.6676161521996592,3,1,26.83081766666667,80.49245300000001,0.5916656861001716,0.5916656861001716,37.49448467092434,337.450362038319,0.3950055687427679,0.3950055687427679),
    @("ECs","Sema3a","Nrp2","FAPs",2,0.6666666666666666,1.397441,4.192323,0.6676161521996591,0.6676161521996592,3,1,10.21969166666667,30.659075,0.2253617819930474,0.2253617819930474,14.28141614235834,128.532745281225,0.1504551657470567,0.1504551657470568),
    @("ECs","Sema3a","Nrp2","sCs",2,0.6666666666666666,1.397441,4.192323,0.6676161521996591,0.6676161521996592,3,1,8.297426666666667,24.89228,0.1829725319067811,0.1829725319067811,11.59516421849333,104.35647796644,0.1221554177098345,0.1221554177098345),
    @("FAPs","Sema3a","Nrp2","ECs",3,1,0.2347746666666667,0.7043240000000001,0.112161700990566,0.112161700990566,3,1,26.83081766666667,80.49245300000001,0.5916656861001716,0.5916656861001716,6.299196274085779,56.69276646677201,0.06636222977074553,0.06636222977074555),
    @("FAPs","Sema3a","Nrp2","FAPs",3,1,0.2347746666666667,0.7043240000000001,0.112161700990566,0.112161700990566,3,1,10.21969166666667,30.659075,0.2253617819930474,0.2253617819930474,2.399324704477778,21.5939223403,0.02527696080660531,0.02527696080660531),
    @("FAPs","Sema3a","Nrp2","sCs",3,1,0.2347746666666667,0.7043240000000001,0.112161700990566,0.112161700990566,3,1,8.297426666666667,24.89228,0.1829725319067811,0.1829725319067811,1.948025579857778,17.53223021872,0.02052251041321518,0.02052251041321518),
    @("sCs","Sema3a","Nrp2","ECs",3,1,0.4609646666666667,1.382894,0.2202221468097748,0.2202221468097748,3,1,26.83081766666667,80.49245300000001,0.5916656861001716,0.5916656861001716,12.36805892210911,111.312530298982,0.1302978875866581,0.1302978875866581),
    @("sCs","Sema3a","Nrp2","FAPs",3,1,0.4609646666666667,1.382894,0.2202221468097748,0.2202221468097748,3,1,10.21969166666667,30.659075,0.2253617819930474,0.2253617819930474,4.710916762561112,42.39825086305,0.04962965543938534,0.04962965543938535),
    @("sCs","Sema3a","Nrp2","sCs",3,1,0.4609646666666667,1.382894,0.2202221468097748,0.2202221468097748,3,1,8.297426666666667,24.89228,0.1829725319067811,0.1829725319067811,3.824820517591111,34.42338465832,0.04029460378373134,0.04029460378373134)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}

# The new table only has 9 data rows (2..10); previously there were 6 (2..7).
# Nothing further to clear since we're growing the range, not shrinking it.
